$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (row 11) to the BZG Departures table, reusing shared
# strings where the values already exist and introducing new ones for the
# values that are new (aircraft model, aircraft id, status time, difference).

$ws.Cells.Item(11, 1).Value = 10                          # NUMBER
$ws.Cells.Item(11, 2).Value = "Monday, Jan 09"             # DATE
$ws.Cells.Item(11, 3).Value = "3:00 PM"                    # TIME
$ws.Cells.Item(11, 4).Value = "LO3994"                     # FLIGHT
$ws.Cells.Item(11, 5).Value = "Warsaw"                     # TO
$ws.Cells.Item(11, 6).Value = "(WAW)"                      # SHORT
$ws.Cells.Item(11, 7).Value = "LOT "                       # AIRLINE
$ws.Cells.Item(11, 8).Value = "E170"                       # MODEL (new)
$ws.Cells.Item(11, 9).Value = "(SP-LDF)"                   # AIRCFAT ID (new)
$ws.Cells.Item(11, 10).Value = "3:15 PM"                   # STATUS (new)
$ws.Cells.Item(11, 11).Borders.LineStyle = 0                # K11 stays blank
$ws.Cells.Item(11, 12).Value = "0 hours, 15 minutes"       # DIFFERENCE (new)
$ws.Cells.Item(11, 13).Borders.LineStyle = 0                # M11 stays blank
